$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "https://kardan.edu.af"
$ws.Range("B9").Value = "NO"
